$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Sending cluster (col A) for the existing data rows switches from MuSCs to Resolving-Mac ----
$ws.Range("A2").Value() = "Resolving-Mac"
$ws.Range("A3").Value() = "Resolving-Mac"
$ws.Range("A4").Value() = "Resolving-Mac"
$ws.Range("A5").Value() = "Resolving-Mac"

# ---- Row 2 (ECs target) : refreshed TPM-derived values ----
$ws.Range("G2").Value() = 0.02530666666666667
$ws.Range("H2").Value() = 0.07592
$ws.Range("M2").Value() = 0.903185
$ws.Range("N2").Value() = 2.709555
$ws.Range("O2").Value() = 0.03154869388788047
$ws.Range("P2").Value() = 0.03154869388788046
$ws.Range("Q2").Value() = 0.02285660173333334
$ws.Range("R2").Value() = 0.2057094156
$ws.Range("S2").Value() = 0.03154869388788047
$ws.Range("T2").Value() = 0.03154869388788046

# ---- Row 3 (FAPs target) : refreshed TPM-derived values ----
$ws.Range("G3").Value() = 0.02530666666666667
$ws.Range("H3").Value() = 0.07592
$ws.Range("O3").Value() = 0.003855347953955327
$ws.Range("P3").Value() = 0.003855347953955326
$ws.Range("Q3").Value() = 0.002793147413333333
$ws.Range("R3").Value() = 0.02513832672
$ws.Range("S3").Value() = 0.003855347953955327
$ws.Range("T3").Value() = 0.003855347953955326

# ---- Row 4 : target cluster switches from MuSCs to Inflammatory-Mac, values refreshed ----
$ws.Range("D4").Value() = "Inflammatory-Mac"
$ws.Range("G4").Value() = 0.02530666666666667
$ws.Range("H4").Value() = 0.07592
$ws.Range("K4").Value() = 1
$ws.Range("L4").Value() = 0.3333333333333333
$ws.Range("M4").Value() = 0.03076233333333333
$ws.Range("N4").Value() = 0.09228699999999999
$ws.Range("O4").Value() = 0.001074543352259254
$ws.Range("P4").Value() = 0.001074543352259254
$ws.Range("Q4").Value() = 0.0007784921155555556
$ws.Range("R4").Value() = 0.00700642904
$ws.Range("S4").Value() = 0.001074543352259254
$ws.Range("T4").Value() = 0.001074543352259254

# ---- Row 5 : target cluster switches from Resolving-Mac to MuSCs, values refreshed ----
$ws.Range("D5").Value() = "MuSCs"
$ws.Range("G5").Value() = 0.02530666666666667
$ws.Range("H5").Value() = 0.07592
$ws.Range("K5").Value() = 3
$ws.Range("L5").Value() = 1
$ws.Range("M5").Value() = 27.52907633333334
$ws.Range("N5").Value() = 82.58722900000001
$ws.Range("O5").Value() = 0.9616041035407232
$ws.Range("P5").Value() = 0.9616041035407231
$ws.Range("Q5").Value() = 0.696669158408889
$ws.Range("R5").Value() = 6.270022425680001
$ws.Range("S5").Value() = 0.9616041035407232
$ws.Range("T5").Value() = 0.9616041035407231

# ---- Row 6 (new) : Resolving-Mac sending -> Resolving-Mac target ----
$ws.Range("A6").Value() = "Resolving-Mac"
$ws.Range("B6").Value() = "Fgf8"
$ws.Range("C6").Value() = "Fgfr4"
$ws.Range("D6").Value() = "Resolving-Mac"
$ws.Range("E6").Value() = 1
$ws.Range("F6").Value() = 0.3333333333333333
$ws.Range("G6").Value() = 0.02530666666666667
$ws.Range("H6").Value() = 0.07592
$ws.Range("I6").Value() = 1
$ws.Range("J6").Value() = 1
$ws.Range("K6").Value() = 2
$ws.Range("L6").Value() = 0.6666666666666666
$ws.Range("M6").Value() = 0.05488933333333334
$ws.Range("N6").Value() = 0.164668
$ws.Range("O6").Value() = 0.001917311265181737
$ws.Range("P6").Value() = 0.001917311265181736
$ws.Range("Q6").Value() = 0.001389066062222222
$ws.Range("R6").Value() = 0.01250159456
$ws.Range("S6").Value() = 0.001917311265181737
$ws.Range("T6").Value() = 0.001917311265181736
